$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - rows 2..25 hold data, F column updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1607
$ws1.Range("F3").Value = 8985
$ws1.Range("F4").Value = 103
$ws1.Range("F6").Value = 687
$ws1.Range("F7").Value = 344
$ws1.Range("F8").Value = 172
$ws1.Range("F9").Value = 48
$ws1.Range("F10").Value = 77
$ws1.Range("F11").Value = 3831
$ws1.Range("F15").Value = 4259
$ws1.Range("F16").Value = 7
$ws1.Range("F20").Value = 331
$ws1.Range("F21").Value = 2
$ws1.Range("F23").Value = 11
$ws1.Range("F24").Value = 2631
$ws1.Range("F25").Value = 111

# Sheet "全部类型" (sheet4) - rows 2..26 hold data (one extra row vs sheet1), F column updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1607
$ws4.Range("F3").Value = 8985
$ws4.Range("F4").Value = 103
$ws4.Range("F6").Value = 687
$ws4.Range("F7").Value = 344
$ws4.Range("F8").Value = 172
$ws4.Range("F9").Value = 48
$ws4.Range("F10").Value = 77
$ws4.Range("F11").Value = 3831
$ws4.Range("F15").Value = 4259
$ws4.Range("F16").Value = 7
$ws4.Range("F20").Value = 331
$ws4.Range("F21").Value = 2
$ws4.Range("F23").Value = 11
$ws4.Range("F24").Value = 2631
$ws4.Range("F26").Value = 111

$wb.Save()
